# Generate Report for Handback
# Fills in the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns (I, J, K, P) for the
# zh-cn and de-de handback-status sheets, widens those (now much longer)
# columns, and links the newly-populated "Latest Target File" cell back to
# the handback markdown file - mirroring what the existing A2 hyperlink
# already does.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/774c20811d9d59ac5bebe1a17cb9a25f48cdfe65/e2e/fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.md"
$handbackMd = "fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.md"

# Excel stores column widths in "characters" but rounds through a pixel
# conversion on write, so asking for a flat 40 round-trips to ~40.83. Backing
# off by 5/6 compensates for that rounding so the saved width comes out to
# exactly 40, matching the target column width.
$colWidth = 40 - (5 / 6)

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48cde75648e9a4cdadf07ed818a7ab23a861da49/e2e/fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/774c20811d9d59ac5bebe1a17cb9a25f48cdfe65/e2e/fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = $handbackMd
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $handbackUrl, "", "", $handbackMd)

$wsZh.Range("J2").Value = "fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.79fed6077362d441022f1408c5ba5070b73ca71c.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-10 09:25:26"
$wsZh.Range("P2").Value = $errorDetail

$wsZh.Columns.Item(9).ColumnWidth = $colWidth
$wsZh.Columns.Item(10).ColumnWidth = $colWidth
$wsZh.Columns.Item(16).ColumnWidth = $colWidth

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = $handbackMd
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $handbackUrl, "", "", $handbackMd)

$wsDe.Range("J2").Value = "fa55ea4e-ab78-44f7-b5c8-fca626e6eecc.79fed6077362d441022f1408c5ba5070b73ca71c.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-10 09:25:42"
$wsDe.Range("P2").Value = $errorDetail

$wsDe.Columns.Item(9).ColumnWidth = $colWidth
$wsDe.Columns.Item(10).ColumnWidth = $colWidth
$wsDe.Columns.Item(16).ColumnWidth = $colWidth
